# Updates cryptos list data (Price column D and Volume(1h) column E)
# Also swaps rows 12 and 13 (TRON / WrappedEther ranking order changed)
# Column D values must stay as text, so force text NumberFormat before
# assigning them (otherwise Excel auto-converts plain-decimal-looking
# strings into numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37's Price (D) value is not changed by this update, so it is
# excluded from the text-format range below (D2:D36 and D38:D51 only).
$ws.Range("D2:D36").NumberFormat = "@"
$ws.Range("D38:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.066.29"
$ws.Range("E2").Value = "  -2.56%  "

$ws.Range("D3").Value = "1.866.90"
$ws.Range("E3").Value = "  -2.02%  "

$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "306.27"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "0.5138"
$ws.Range("E7").Value = "  -1.84%  "

$ws.Range("D8").Value = "0.3756"
$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("D9").Value = "0.07162"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "0.8892"
$ws.Range("E10").Value = "  -1.96%  "

$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  -2.95%  "

# Row 12 and 13 swap: WrappedEther moves to row 12, TRON moves to row 13
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.887.74"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07608"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "89.66"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "0.000008471"
$ws.Range("E17").Value = "  -2.58%  "

$ws.Range("D18").Value = "14.08"
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").Value = "27.091.39"
$ws.Range("E20").Value = "  -2.57%  "

$ws.Range("D21").Value = "5.034"
$ws.Range("E21").Value = "  -2.17%  "

$ws.Range("D22").Value = "2.114.64"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").Value = "10.50"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("D24").Value = "6.462"
$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D25").Value = "1.839"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "147.62"
$ws.Range("E26").Value = "  -4.01%  "

$ws.Range("D27").Value = "17.98"
$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("D28").Value = "2.110"
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "112.75"
$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("D30").Value = "4.658"
$ws.Range("E30").Value = "  -3.87%  "

$ws.Range("D31").Value = "4.703"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("D32").Value = "0.09099"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").Value = "0.05124"
$ws.Range("E33").Value = "  -2.89%  "

$ws.Range("D34").Value = "3.070"
$ws.Range("E34").Value = "  -3.29%  "

$ws.Range("D35").Value = "1.158"
$ws.Range("E35").Value = "  -5.74%  "

$ws.Range("D36").Value = "0.7249"
$ws.Range("E36").Value = "  -7.05%  "

$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").Value = "2.496"
$ws.Range("E38").Value = "  -4.73%  "

$ws.Range("D39").Value = "3.041"
$ws.Range("E39").Value = "  -1.09%  "

$ws.Range("D40").Value = "1.073"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("D41").Value = "0.5337"
$ws.Range("E41").Value = "  -3.58%  "

$ws.Range("D42").Value = "6.559"
$ws.Range("E42").Value = "  -1.82%  "

$ws.Range("D43").Value = "115.85"
$ws.Range("E43").Value = "  +1.19%  "

$ws.Range("D44").Value = "8.298"
$ws.Range("E44").Value = "  -3.05%  "

$ws.Range("D45").Value = "0.1467"
$ws.Range("E45").Value = "  -2.93%  "

$ws.Range("D46").Value = "0.4640"
$ws.Range("E46").Value = "  -3.38%  "

$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").Value = "9.978"
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("D49").Value = "1.573"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("D50").Value = "36.55"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").Value = "63.98"
$ws.Range("E51").Value = "  -4.31%  "
